$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings (not numeric values, e.g. "24.151.33"),
# so force the column to Text format to avoid Excel auto-converting
# values that look numeric (e.g. "0.9997", "158.00", "1.0000").
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "24.151.33"
$ws.Range("E2").Value = "  -2.84%  "
$ws.Range("D3").Value = "1.642.91"
$ws.Range("E3").Value = "  -2.65%  "
$ws.Range("D4").Value = "0.9997"
$ws.Range("E4").Value = "  -0.48%  "
$ws.Range("D5").Value = "308.77"
$ws.Range("E5").Value = "  -1.90%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.36%  "
$ws.Range("E7").Value = "  -1.12%  "
$ws.Range("D8").Value = "0.3865"
$ws.Range("E8").Value = "  -3.09%  "
$ws.Range("D9").Value = "1.003"
$ws.Range("E9").Value = "  -0.21%  "
$ws.Range("D10").Value = "49.62"
$ws.Range("E10").Value = "  -4.69%  "
$ws.Range("D11").Value = "1.350"
$ws.Range("E11").Value = "  -5.99%  "
$ws.Range("D12").Value = "0.08660"
$ws.Range("E12").Value = "  -0.58%  "
$ws.Range("E13").Value = "  -6.83%  "
$ws.Range("D14").Value = "7.099"
$ws.Range("E14").Value = "  -3.11%  "
$ws.Range("D15").Value = "0.00001290"
$ws.Range("E15").Value = "  -2.86%  "
$ws.Range("D16").Value = "7.459"
$ws.Range("E16").Value = "  -4.89%  "
$ws.Range("D17").Value = "1.642.23"
$ws.Range("E17").Value = "  -1.29%  "
$ws.Range("D18").Value = "95.29"
$ws.Range("E18").Value = "  +1.01%  "
$ws.Range("D19").Value = "0.06916"
$ws.Range("E19").Value = "  -3.19%  "
$ws.Range("D20").Value = "20.57"
$ws.Range("E20").Value = "  +1.55%  "
$ws.Range("D21").Value = "6.901"
$ws.Range("E21").Value = "  -4.02%  "
$ws.Range("D22").Value = "0.9996"
$ws.Range("E22").Value = "  -0.53%  "
$ws.Range("D23").Value = "13.57"
$ws.Range("E23").Value = "  -4.23%  "
$ws.Range("D24").Value = "24.140.48"
$ws.Range("E24").Value = "  -2.92%  "
$ws.Range("D25").Value = "2.330"
$ws.Range("E25").Value = "  -2.67%  "
$ws.Range("D26").Value = "2.764"
$ws.Range("E26").Value = "  -1.88%  "
$ws.Range("D27").Value = "22.35"
$ws.Range("E27").Value = "  -2.83%  "
$ws.Range("D28").Value = "158.00"
$ws.Range("E28").Value = "  -2.45%  "
$ws.Range("D29").Value = "8.487"
$ws.Range("E29").Value = "  +6.82%  "
$ws.Range("D30").Value = "140.18"
$ws.Range("E30").Value = "  -5.94%  "
$ws.Range("D31").Value = "5.303"
$ws.Range("E31").Value = "  -11.28%  "
$ws.Range("D32").Value = "2.421"
$ws.Range("E32").Value = "  -8.70%  "
$ws.Range("D33").Value = "1.810.41"
$ws.Range("E33").Value = "  -1.98%  "
$ws.Range("D34").Value = "6.918"
$ws.Range("E34").Value = "  -1.12%  "
$ws.Range("D35").Value = "0.08028"
$ws.Range("E35").Value = "  -5.26%  "
$ws.Range("D36").Value = "0.02897"
$ws.Range("E36").Value = "  -6.22%  "
$ws.Range("D37").Value = "0.2678"
$ws.Range("E37").Value = "  -5.99%  "
$ws.Range("D38").Value = "0.9463"
$ws.Range("E38").Value = "  -7.61%  "
$ws.Range("D39").Value = "0.09206"
$ws.Range("E39").Value = "  -4.11%  "
$ws.Range("D40").Value = "1.470"
$ws.Range("E40").Value = "  +0.60%  "
$ws.Range("D41").Value = "9.920"
$ws.Range("E41").Value = "  -6.97%  "
$ws.Range("D42").Value = "0.7551"
$ws.Range("E42").Value = "  -6.24%  "
$ws.Range("E43").Value = "  -5.64%  "
$ws.Range("D44").Value = "15.94"
$ws.Range("E44").Value = "  -5.12%  "
$ws.Range("D45").Value = "0.6908"
$ws.Range("E45").Value = "  -4.55%  "
$ws.Range("D46").Value = "2.463"
$ws.Range("E46").Value = "  -5.79%  "
$ws.Range("D47").Value = "4.090"
$ws.Range("E47").Value = "  -2.70%  "
$ws.Range("D48").Value = "1.0000"
$ws.Range("E48").Value = "  -0.33%  "
$ws.Range("E49").Value = "  -5.09%  "
$ws.Range("D50").Value = "1.258"
$ws.Range("E50").Value = "  -8.49%  "
$ws.Range("D51").Value = "133.06"
$ws.Range("E51").Value = "  -4.22%  "
